# Update the "High School Year attended" label/example cell and its
# corresponding answer cell on Sheet1, then leave the selection on A8
# (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B8 first so the new shared-string entries are appended in the same
# order as the target workbook ("Junior" before the updated label text).
$ws.Range("B8").Value = "Junior"
$ws.Range("A8").Value = "High School Year attended: (e.g. Sophomore)"

$ws.Range("A8").Select()
